# Actualización desde MV -datos-
# Append 4 new daily rows (08-10-2021, 12-10-2021, 13-10-2021, 14-10-2021)
# to the bottom of the "Tasas de captación" table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date labels must land as plain text (shared strings), matching the
# existing "Serie" column — not get auto-converted to Excel date serials.
# Briefly forcing Text format before the assignment, then clearing the
# format again afterwards, keeps the cell on the default (unstyled) xf
# while still storing the value as text.
$dates = @("08-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

$data = @(
    @(0.16, 0.23, 0.01, 0.18, 0.15),
    @(0.17, 0.21, 0.02, 0.18, 0.38),
    @(0.18, 0.28, 0.01, 0.26, 0.3),
    @(0.21, 0.3, 0.07000000000000001, 0.22, 0.52)
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 195 + $i

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dates[$i]
    $dateCell.ClearFormats()

    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $values[$c]
    }
}
